# C5-PowerPoint.pptx edit
#
# 1) The table on slide 6 gets a new built-in table style
#    ({BC28380E-B53B-45AE-81C8-07B3D8DD8FDC} -> {B103EAF7-A487-45D2-8CD8-EEE3CF706473}).
#
# 2) The design's theme colour palette is switched from the custom "Integral"
#    palette over to the stock Office palette (the deck's two embedded theme
#    parts trade places: the slide master/point-of-use theme ends up holding
#    the default "Office Theme" colours). We reproduce the colour values of
#    that swap through the PowerPoint object model's theme colour scheme,
#    which is the surface that actually persists into the theme part.

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 6 --------------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{B103EAF7-A487-45D2-8CD8-EEE3CF706473}", $true)

# --- 2) Theme colour palette -----------------------------------------------
# Target palette ("Office Theme"): dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeRGB = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeRGB[$i - 1]
}
